$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "val3"
$ws.Range("B4").Value = 44185
